$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L2").Value = 0.7414799316992706
$ws.Range("Q2").Value = -0.008954984241185193

$ws.Range("L3").Value = 0.7414799316992706
$ws.Range("Q3").Value = -0.008954984241185193

$ws.Range("L4").Value = 0.7088701277621434
$ws.Range("Q4").Value = -0.016229058904976

$ws.Range("L5").Value = 0.7808705382933501
$ws.Range("Q5").Value = 0.02351246133036377

$ws.Range("L6").Value = 0.7808705382933501
$ws.Range("Q6").Value = 0.02351246133036377

$ws.Range("L7").Value = 0.7832122299475502
$ws.Range("Q7").Value = -0.01081937260331701
